$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows above the old row 10 (pushes old rows 10:13 down to 17:20)
$ws.Rows("10:16").Insert()

# New "alfa" iteration-error summary block (rows 10-12)
$ws.Range("B10").Value = "alfa 1"
$ws.Range("C10").Formula = "=+ABS(C6)+ABS(D6)"

$ws.Range("B11").Value = "alfa 2"
$ws.Range("C11").Formula = "=+ABS(B7)+ABS(D7)"

# New highlighted "alfa" conclusion row (row 14) -- written before "alfa 3" below
# so the shared-string table gets the same insertion order as the source file
$ws.Range("B14").Value = "alfa"
$ws.Range("C14").Value = 0.1

$ws.Range("B12").Value = "alfa 3"
$ws.Range("C12").Formula = "=+ABS(B8)+ABS(C8)"

$ws.Range("D14").Value = " menor a 1 por lo tanto covergera rapido"

$ws.Range("B14:C14").Interior.Color = 15773696

# Update the view to match the saved state (scroll so row 4 is at the top,
# and leave the selection on G14)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("G14").Select()
